# Auto update Excel log
# Appends new sensor-log rows to the mmWave, PIR and Humidity sheets.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$Row,
        [string]$Date,
        [string]$Timestamp,
        [string]$Hour,
        [string]$Location,
        [string]$Value,
        [string]$Status
    )

    # Force text format so Excel does not reinterpret dates/times/percentages
    # as numeric values.
    $ws.Cells.Item($Row, 1).NumberFormat = "@"
    $ws.Cells.Item($Row, 1).Value = $Date

    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 2).Value = $Timestamp

    $ws.Cells.Item($Row, 3).NumberFormat = "@"
    $ws.Cells.Item($Row, 3).Value = $Hour

    $ws.Cells.Item($Row, 4).NumberFormat = "@"
    $ws.Cells.Item($Row, 4).Value = $Location

    $ws.Cells.Item($Row, 5).NumberFormat = "@"
    $ws.Cells.Item($Row, 5).Value = $Value

    $ws.Cells.Item($Row, 6).NumberFormat = "@"
    $ws.Cells.Item($Row, 6).Value = $Status
}

# ---------------------------------------------------------------------------
# mmWave sheet - add rows 23..26 (Living Room presence log, ending in a fall)
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

Add-LogRow $wsMmWave 23 "2026-01-30" "14:53:18" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 24 "2026-01-30" "14:53:28" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 25 "2026-01-30" "14:53:39" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 26 "2026-01-30" "14:54:01" "14:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"

# ---------------------------------------------------------------------------
# PIR sheet - add rows 7..13 (Bathroom no-motion log)
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

Add-LogRow $wsPIR 7  "2026-01-30" "14:53:16" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 8  "2026-01-30" "14:53:21" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 9  "2026-01-30" "14:53:26" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 10 "2026-01-30" "14:53:31" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 11 "2026-01-30" "14:53:36" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 12 "2026-01-30" "14:53:41" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 13 "2026-01-30" "14:53:46" "14:00" "Bathroom" "No Motion" "Inactive"

# ---------------------------------------------------------------------------
# Humidity sheet - add rows 6..12 (Bathroom humidity readings)
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

Add-LogRow $wsHumidity 6  "2026-01-30" "14:53:17" "14:00" "Bathroom" "88.4%" "Active"
Add-LogRow $wsHumidity 7  "2026-01-30" "14:53:22" "14:00" "Bathroom" "88.4%" "Active"
Add-LogRow $wsHumidity 8  "2026-01-30" "14:53:27" "14:00" "Bathroom" "88.4%" "Active"
Add-LogRow $wsHumidity 9  "2026-01-30" "14:53:32" "14:00" "Bathroom" "88.3%" "Active"
Add-LogRow $wsHumidity 10 "2026-01-30" "14:53:37" "14:00" "Bathroom" "87.4%" "Active"
Add-LogRow $wsHumidity 11 "2026-01-30" "14:53:42" "14:00" "Bathroom" "88.3%" "Active"
Add-LogRow $wsHumidity 12 "2026-01-30" "14:53:47" "14:00" "Bathroom" "88.3%" "Active"
